$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

Set-TextValue 'D2' '30.779.79'
Set-TextValue 'E2' '  +0.74%  '
Set-TextValue 'D3' '1.895.55'
Set-TextValue 'E3' '  +1.18%  '
Set-TextValue 'D4' '0.9999'
Set-TextValue 'E4' '  +0.04%  '
Set-TextValue 'D5' '247.13'
Set-TextValue 'E5' '  +0.05%  '
Set-TextValue 'D6' '0.9996'
Set-TextValue 'E6' '  +0.03%  '
Set-TextValue 'D7' '0.4735'
Set-TextValue 'E7' '  +0.10%  '
Set-TextValue 'D8' '0.2931'
Set-TextValue 'E8' '  +0.42%  '
Set-TextValue 'D9' '0.06520'
Set-TextValue 'E9' '  +0.39%  '
Set-TextValue 'D10' '22.66'
Set-TextValue 'E10' '  +2.01%  '
Set-TextValue 'D11' '0.07793'
Set-TextValue 'E11' '  +1.01%  '
Set-TextValue 'D12' '0.7421'
Set-TextValue 'E12' '  +0.01%  '
Set-TextValue 'B13' 'Litecoin'
Set-TextValue 'C13' 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
Set-TextValue 'D13' '97.01'
Set-TextValue 'E13' '  -0.47%  '
Set-TextValue 'B14' 'WrappedEther'
Set-TextValue 'C14' 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
Set-TextValue 'D14' '1.887.89'
Set-TextValue 'E14' '  +0.97%  '
Set-TextValue 'D15' '5.248'
Set-TextValue 'E15' '  +1.70%  '
Set-TextValue 'D16' '286.15'
Set-TextValue 'E16' '  +4.24%  '
Set-TextValue 'D17' '30.755.56'
Set-TextValue 'E17' '  +0.75%  '
Set-TextValue 'D18' '13.28'
Set-TextValue 'E18' '  -0.88%  '
Set-TextValue 'D19' '0.000007535'
Set-TextValue 'E19' '  +0.20%  '
Set-TextValue 'E20' '  -0.03%  '
Set-TextValue 'D21' '2.131.36'
Set-TextValue 'E21' '  +0.81%  '
Set-TextValue 'D22' '5.339'
Set-TextValue 'E22' '  +1.41%  '
Set-TextValue 'D23' '0.9999'
Set-TextValue 'E23' '  +0.02%  '
Set-TextValue 'D24' '6.291'
Set-TextValue 'E24' '  +1.88%  '
Set-TextValue 'D25' '9.230'
Set-TextValue 'E25' '  -0.58%  '
Set-TextValue 'D26' '164.92'
Set-TextValue 'E26' '  +0.51%  '
Set-TextValue 'D27' '19.04'
Set-TextValue 'E27' '  +1.41%  '
Set-TextValue 'D28' '1.927'
Set-TextValue 'E28' '  +0.37%  '
Set-TextValue 'E29' '  -1.84%  '
Set-TextValue 'D30' '1.344'
Set-TextValue 'E30' '  -0.46%  '
Set-TextValue 'D31' '1.491'
Set-TextValue 'E31' '  -1.24%  '
Set-TextValue 'D32' '4.312'
Set-TextValue 'E32' '  +0.44%  '
Set-TextValue 'D33' '4.189'
Set-TextValue 'E33' '  +1.84%  '
Set-TextValue 'D34' '0.04910'
Set-TextValue 'E34' '  +2.06%  '
Set-TextValue 'D35' '1.136'
Set-TextValue 'E35' '  +1.47%  '
Set-TextValue 'D36' '0.7001'
Set-TextValue 'E36' '  +0.58%  '
Set-TextValue 'D37' '2.709'
Set-TextValue 'E37' '  -0.07%  '
Set-TextValue 'D38' '0.01904'
Set-TextValue 'E38' '  +2.73%  '
Set-TextValue 'D39' '2.840'
Set-TextValue 'E39' '  +3.49%  '
Set-TextValue 'D40' '76.14'
Set-TextValue 'E40' '  +3.36%  '
Set-TextValue 'D41' '6.325'
Set-TextValue 'E41' '  +1.50%  '
Set-TextValue 'D42' '2.017'
Set-TextValue 'E42' '  +2.52%  '
Set-TextValue 'D43' '0.4306'
Set-TextValue 'E43' '  +2.75%  '
Set-TextValue 'D44' '0.9999'
Set-TextValue 'E44' '  +0.06%  '
Set-TextValue 'D45' '0.8372'
Set-TextValue 'E45' '  +0.37%  '
Set-TextValue 'D46' '101.88'
Set-TextValue 'E46' '  -0.73%  '
Set-TextValue 'D47' '9.634'
Set-TextValue 'E47' '  +2.94%  '
Set-TextValue 'D48' '7.047'
Set-TextValue 'E48' '  +0.92%  '
Set-TextValue 'D49' '35.45'
Set-TextValue 'E49' '  +0.06%  '
Set-TextValue 'D50' '912.33'
Set-TextValue 'E50' '  -1.61%  '
Set-TextValue 'D51' '0.3983'
Set-TextValue 'E51' '  +2.82%  '
